# Insert a blank row above row 13 ("Length of stay" header row).
# This pushes the whole "Length of stay" table (previously rows 13-24)
# down by one row (to rows 14-25), extending the used range to A1:L25.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(13).Insert()

# Match the resulting selection left behind by the edit.
$ws.Range("A14:L25").Select() | Out-Null
